$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) columns for rows with numeric-only changes.
# Some price values are numeric-looking strings (e.g. "1.0000") that must stay as
# TEXT (matching the source inlineStr cells) rather than being auto-converted to
# numbers by Excel. We force text via NumberFormat="@", assign the value, then
# ClearFormats() so the cell keeps no explicit style (matching the original, which
# has no "s" attribute on these cells).

$ws.Range("D2").Value = "30.178.75"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "1.862.97"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4668"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2827"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06529"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.82"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07860"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.12"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.09%  "
$ws.Range("D13").Value = "1.863.11"
$ws.Range("E13").Value = "  -2.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.113"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6698"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "278.06"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.57%  "
$ws.Range("D17").Value = "30.176.03"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9992"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.460"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.61"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").Value = "2.102.18"
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007236"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.76%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.144"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.288"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.03"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.94"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.896"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -7.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.353"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09569"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.388"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.468"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.096"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04679"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6997"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.096"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.717"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01855"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.340"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -6.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.525"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.73"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.928"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8447"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.36%  "

# Rows 44 and 45: coin order swapped (PaxDollar <-> TheSandbox) with updated price/volume
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4159"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.87%  "

# Continue updating price and volume columns (D, E) for remaining rows
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.84"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "958.43"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.130"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.240"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.93"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1131"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.62%  "
